$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = 101
$ws.Range("B12").Value = "Hidden"
$ws.Range("C12").Value = "A"
$ws.Range("D12").Value = "Laos"

$ws.Range("A13").Value = 102
$ws.Range("B13").Value = "Hidden"
$ws.Range("C13").Value = "B"
$ws.Range("D13").Value = "Italy"

$ws.Rows.Item(12).Hidden = $true
$ws.Rows.Item(13).Hidden = $true

[void]$ws.Range("H22").Select()
